# Update the agent_size workbook:
#  - voxel_total (column B) values shrink for every space on Sheet1
#  - voxel_depth (column D) drops from 7 to 4 for every space
#  - the sheet view zoom level / selection reflect the author's new position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# voxel_total (column B) new values, rows 2-20 (student_housing ... arcade)
$voxelTotals = @(91, 80, 340, 20, 23, 6, 6, 14, 8, 26, 11, 57, 14, 17, 1, 40, 1, 28, 11)

for ($i = 0; $i -lt $voxelTotals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $voxelTotals[$i]
}

# voxel_depth (column D) goes from 7 to 4 for every data row
$ws.Range("D2:D20").Value = 4

# View state: zoom in a bit and move the selection down one row
$excel.ActiveWindow.Zoom = 84
$ws.Range("B21").Select()
